$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Per Package" / "Per Handling Unit" header cells to the new
# terminology used by the updated manifest/reference workflow.
$ws.Range("J2").Value = "PCS per PU"
$ws.Range("K2").Value = "PU per HU"

# Move the active selection to K3 (matches the saved cursor position).
$ws.Range("K3").Select() | Out-Null
